$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (new quarter columns), shifting
# the existing quarterly data from D:K to F:M.
$ws.Range("D:E").Insert()

# The inserted columns carry no formatting; copy number formats from the
# first surviving data column (F, the old D) across the same rows so the
# new cells render as dates / thousands-formatted numbers like their
# neighbours.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Range("E7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Rows 36/37 and 78/79 are spacer rows with no data in columns D:K in the
# original sheet; undo the incidental formatting the bulk paste above
# applied to their D:E cells so they stay empty, matching the source rows.
$ws.Range("D36:E36").Clear()
$ws.Range("D37:E37").Clear()
$ws.Range("D78:E78").Clear()
$ws.Range("D79:E79").Clear()

# Populate the two new quarter columns (D = Q2'19, E = Q1'19) with the
# latest reported figures for every statement line.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 344800
$ws.Range("E8").Value = 348100
$ws.Range("D9").Value = 154700
$ws.Range("E9").Value = 152900
$ws.Range("D10").Value = 190100
$ws.Range("E10").Value = 195200
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 200
$ws.Range("E14").Value = 100
$ws.Range("D15").Value = 103900
$ws.Range("E15").Value = 109100
$ws.Range("D17").Value = 340100
$ws.Range("E17").Value = 347700
$ws.Range("D18").Value = 4700
$ws.Range("E18").Value = 400
$ws.Range("D20").Value = -25600
$ws.Range("E20").Value = -24200
$ws.Range("D21").Value = -26100
$ws.Range("E21").Value = -26400
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = -20900
$ws.Range("E23").Value = -23800
$ws.Range("D24").Value = -6100
$ws.Range("E24").Value = -4600
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -14800
$ws.Range("E26").Value = -19200
$ws.Range("D27").Value = -14900
$ws.Range("E27").Value = -19500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 800
$ws.Range("E29").Value = 4400
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 25600
$ws.Range("E32").Value = 24200
$ws.Range("D33").Value = -14100
$ws.Range("E33").Value = -15100
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -14100
$ws.Range("E35").Value = -15100
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 9600
$ws.Range("E41").Value = 3800
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 144200
$ws.Range("E43").Value = 155500
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 44300
$ws.Range("E45").Value = 40600
$ws.Range("D46").Value = 198100
$ws.Range("E46").Value = 199900
$ws.Range("D47").Value = 110900
$ws.Range("E47").Value = 110700
$ws.Range("D48").Value = 1927100
$ws.Range("E48").Value = 1955800
$ws.Range("D49").Value = 1275700
$ws.Range("E49").Value = 1292900
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 23400
$ws.Range("E52").Value = 36700
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 3535300
$ws.Range("E54").Value = 3596000
$ws.Range("D57").Value = 32500
$ws.Range("E57").Value = 15700
$ws.Range("D58").Value = 30500
$ws.Range("E58").Value = 31800
$ws.Range("D59").Value = 220600
$ws.Range("E59").Value = 231000
$ws.Range("D60").Value = 283600
$ws.Range("E60").Value = 278500
$ws.Range("D61").Value = 2303600
$ws.Range("E61").Value = 2302800
$ws.Range("D62").Value = 532400
$ws.Range("E62").Value = 526200
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 3125500
$ws.Range("E66").Value = 3113400
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -50800
$ws.Range("E72").Value = -36900
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 409700
$ws.Range("E76").Value = 482600
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -14100
$ws.Range("E81").Value = -15100
$ws.Range("D83").Value = 0
$ws.Range("E83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 93300
$ws.Range("E89").Value = 69700
$ws.Range("D91").Value = -58100
$ws.Range("E91").Value = -61900
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -57600
$ws.Range("E94").Value = -40700
$ws.Range("D96").Value = -27600
$ws.Range("E96").Value = -27600
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -29900
$ws.Range("E100").Value = -35800
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 5800
$ws.Range("E102").Value = -6800

# A handful of prior-quarter figures (old columns F/G, now H/I) were
# restated along with this update.
$ws.Range("H9").Value = 155500
$ws.Range("I9").Value = 148400
$ws.Range("H10").Value = 200900
$ws.Range("I10").Value = 214900
$ws.Range("H17").Value = 349400
$ws.Range("I17").Value = 371000
$ws.Range("H18").Value = 7000
$ws.Range("I18").Value = -7700
$ws.Range("H20").Value = -22100
$ws.Range("I20").Value = -27000
$ws.Range("H32").Value = 22100
$ws.Range("I32").Value = 27000

Write-Host "Applied quarterly update."
